# CIERRE 26 MAR 22
# Roll the weekly payroll receipt sheet ("Hoja1") from SEMANA 11 (14-20 MAR 2022)
# to SEMANA 12 (21-27 MAR 2022), updating the period label, today's date stamps,
# and the period's pay figures (loan/"PRESTAMO" line replaced by "EXTRA" with a
# value, the social-security withholding reduced, and the second slip's bonus
# moved out of the "# 21" line and into the total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# --- Week label (shared by both slips via formulas off B9) ---------------
$ws.Range("B9").Value2 = "SEMANA   12  DEL    21      Al   27   DE   MARZO          2022"

# --- First slip: social-security figure (K21) drops, and the old
#     "PRESTAMO" (loan) line becomes "EXTRA" with a 500 value ------------
$ws.Range("K21").Value2 = 1120
$ws.Range("D25").Value2 = "EXTRA"
$ws.Range("E25").Value2 = 500

# --- Second slip: the flat 1250 bonus on row 40 goes to 0 (folded into
#     the total instead) ---------------------------------------------------
$ws.Range("K40").Value2 = 0

# --- Selection / view state left by the author on save -------------------
$ws.Range("E54").Select()
